$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: paragraph "Drawing pictures and making all of the visuals to
# the game" - collapse the three runs (with the gramStart/gramEnd
# proofErr markers around "all of") into a single run, and mark the
# start of the run with a <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------
$target = "Drawing pictures and making all of the visuals to the game"

$found = $d.Content.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = $d.Paragraphs.Item(8)
$r = $p.Range

# Insert a new run at the very start of the paragraph that carries the
# lastRenderedPageBreak marker together with the full replacement text.
$rStart = $d.Range($r.Start, $r.Start)
$rStart.InsertXML("<w:p><w:r><w:lastRenderedPageBreak/><w:t>" + $target + "</w:t></w:r></w:p>")

# The original (three-run) text now immediately follows the text we
# just inserted - remove it, leaving only the single merged run behind.
$p2 = $d.Paragraphs.Item(8)
$r2 = $p2.Range
$cutStart = $r2.Start + $target.Length
$cutEnd = $r2.End
$rDel = $d.Range($cutStart, $cutEnd)
$rDel.Delete()

# ---------------------------------------------------------------------
# Edit 2: add a new "Database" bullet right after "Python for Unity ".
# ---------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("Python for Unity", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pyPara = $d.Paragraphs.Item(16)
$pyPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(17)
$newPara.Range.Text = "Database"
